# Refresh "想去人数" (column F) figures across all four sheets,
# matching the data snapshot committed at 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 22
$ws.Cells.Item(5, 6).Value = 6206
$ws.Cells.Item(6, 6).Value = 694
$ws.Cells.Item(7, 6).Value = 1076
$ws.Cells.Item(8, 6).Value = 51
$ws.Cells.Item(9, 6).Value = 99
$ws.Cells.Item(12, 6).Value = 638
$ws.Cells.Item(13, 6).Value = 1099
$ws.Cells.Item(14, 6).Value = 69
$ws.Cells.Item(15, 6).Value = 345
$ws.Cells.Item(17, 6).Value = 6
$ws.Cells.Item(18, 6).Value = 1395
$ws.Cells.Item(19, 6).Value = 629
$ws.Cells.Item(21, 6).Value = 370
$ws.Cells.Item(23, 6).Value = 1046
$ws.Cells.Item(24, 6).Value = 99
$ws.Cells.Item(25, 6).Value = 2128
$ws.Cells.Item(26, 6).Value = 217
$ws.Cells.Item(27, 6).Value = 68
$ws.Cells.Item(28, 6).Value = 377
$ws.Cells.Item(29, 6).Value = 52
$ws.Cells.Item(30, 6).Value = 3445

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(9, 6).Value = 678
$ws.Cells.Item(19, 6).Value = 369
$ws.Cells.Item(20, 6).Value = 309
$ws.Cells.Item(27, 6).Value = 81

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 1166
$ws.Cells.Item(8, 6).Value = 1542
$ws.Cells.Item(12, 6).Value = 711

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 1166
$ws.Cells.Item(6, 6).Value = 1542
$ws.Cells.Item(9, 6).Value = 711
$ws.Cells.Item(13, 6).Value = 6206
$ws.Cells.Item(15, 6).Value = 694
$ws.Cells.Item(16, 6).Value = 1076
$ws.Cells.Item(17, 6).Value = 678
$ws.Cells.Item(18, 6).Value = 51
$ws.Cells.Item(19, 6).Value = 99
$ws.Cells.Item(22, 6).Value = 638
$ws.Cells.Item(27, 6).Value = 1099
$ws.Cells.Item(28, 6).Value = 69
$ws.Cells.Item(29, 6).Value = 345
$ws.Cells.Item(31, 6).Value = 369
$ws.Cells.Item(33, 6).Value = 1395
$ws.Cells.Item(34, 6).Value = 629
$ws.Cells.Item(35, 6).Value = 370
$ws.Cells.Item(39, 6).Value = 1046
$ws.Cells.Item(40, 6).Value = 99
$ws.Cells.Item(42, 6).Value = 2128
$ws.Cells.Item(44, 6).Value = 217
$ws.Cells.Item(45, 6).Value = 68
$ws.Cells.Item(46, 6).Value = 377
$ws.Cells.Item(47, 6).Value = 52
$ws.Cells.Item(48, 6).Value = 3445
